# Adding OPR293 test scenario
# Append 4 new data rows (rows 3-6) to the first worksheet
# (LTE001_ACC_00001 / sheet1.xml), matching the existing table's columns:
# AgentCode, ShipperCode, ConsigneeCode, Origin, Destination, ProductCode,
# SCC, Commodity, ShipmentDescription, ServiceCargoClass, Piece, Weight,
# ChargeType, ModeOfPayment, cartType

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LTE001_ACC_00001")

# Make this the active/selected sheet (it becomes the first/active tab)
$ws.Activate()

$newRows = @(
    @(11377, 11377, 11377, "SEA", "LAX", "GOLDSTREAK",  "None", "NONSCR", "None", "None", 20, 750, "PP", "CREDIT", "CART"),
    @(11377, 11377, 11377, "DFW", "SEA", "GENERAL",     "None", "NONSCR", "None", "None",  1, 100, "PP", "CREDIT", "CART"),
    @(11377, 11377, 11377, "DFW", "SEA", "PRIORITY",    "None", "NONSCR", "None", "None",  8, 600, "PP", "CREDIT", "CART"),
    @(11377, 11377, 11377, "SEA", "ANC", "GENERAL",     "None", "NONSCR", "None", "None",  1, 100, "PP", "CREDIT", "CART")
)

$rowIndex = 3
foreach ($rowValues in $newRows) {
    $colIndex = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

# Match the final selection left behind in the sheet
$ws.Range("F6").Select()
